# Updated symbol list on Sat Dec 24 07:45:28 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking quotes as TEXT (inlineStr in
# the original OOXML), so each price cell is forced to text format ("@")
# before the new value is written -- otherwise the COM layer would coerce
# the numeric-looking string into a real number. The "Volume(1h)" column (E)
# holds free-form text already, so those two cells are updated directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $value) {
    $cell = $ws.Cells.Item($row, 4)   # column D = Price
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-PriceText 2  "245.71"
Set-PriceText 3  "22.10"
Set-PriceText 4  "5.378"
Set-PriceText 5  "0.05973"
Set-PriceText 7  "6.392"
Set-PriceText 8  "0.8105"
Set-PriceText 9  "0.9552"
Set-PriceText 10 "0.1430"
Set-PriceText 11 "0.07395"
Set-PriceText 12 "0.03345"
Set-PriceText 13 "0.03065"
Set-PriceText 14 "0.09409"
Set-PriceText 15 "3.998"
Set-PriceText 16 "0.001606"
Set-PriceText 17 "0.04803"
Set-PriceText 18 "0.0005909"

$ws.Cells.Item(18, 5).Value = "17OneONE"   # column E = Volume(1h)

Set-PriceText 19 "0.006113"
Set-PriceText 20 "0.005115"
Set-PriceText 21 "0.0009850"
Set-PriceText 22 "0.00006999"
Set-PriceText 23 "3.746"
Set-PriceText 24 "2.184"
Set-PriceText 26 "0.1333"
Set-PriceText 27 "0.0002461"
Set-PriceText 40 "0.03995"
Set-PriceText 41 "0.006579"
Set-PriceText 43 "0.002899"
Set-PriceText 44 "0.005835"
Set-PriceText 45 "0.00005254"
Set-PriceText 47 "0.8499"
Set-PriceText 48 "0.03161"

$ws.Cells.Item(48, 5).Value = "47BOLOBOLOWorstin24h"   # column E = Volume(1h)

Set-PriceText 49 "0.00002100"
